# Correccion a Diebold Mariano y revision de Cap1
# Update the DM-test summary table: recompute Proporcion_Sig (B), ECRPS_Mejor (C),
# Mejor_N_Calib (D) and DM statistic (E) for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (AREPD): Proporcion_Sig 4/10 -> 3/10 ; ECRPS_Mejor 102.4 -> 76.8
$ws.Range("B3").Value = "3/10"
$ws.Range("C3").Value = 76.8

# Row 4 (Block Bootstrapping): Proporcion_Sig 4/10 -> 3/10 ; ECRPS_Mejor 102.4 -> 76.8
$ws.Range("B4").Value = "3/10"
$ws.Range("C4").Value = 76.8

# Row 5 (AV-MCPS): Proporcion_Sig 3/10 -> 1/10 ; ECRPS_Mejor 76.8 -> 25.6 ;
#                  Mejor_N_Calib 60 -> 40 ; DM statistic updated
$ws.Range("B5").Value = "1/10"
$ws.Range("C5").Value = 25.6
$ws.Range("D5").Value = 40
$ws.Range("E5").Value = 2.789165874121422

# Row 6 (MCPS): Proporcion_Sig 2/10 -> 0/10 ; ECRPS_Mejor 51.2 -> 0 ;
#               Mejor_N_Calib 40 -> 60 ; DM statistic updated
$ws.Range("B6").Value = "0/10"
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 60
$ws.Range("E6").Value = 2.89488124401165

$wb.Save()
